$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values (columns B,C,D,F,G,L,M,N for rows 2-25)
$data = @{
    2 = @{ "B"="22.59366136756072"; "C"="9.767139270066103"; "D"="7.827313039962065"; "F"="41.16730046323628"; "G"="3.70617917375858"; "L"="10.3155031456119"; "M"="18.25970347627684"; "N"="21.06090081141447" }
    3 = @{ "B"="22.16509135479668"; "C"="9.154281415116335"; "D"="7.849477109082049"; "F"="40.83119117588051"; "G"="3.710603636976363"; "L"="10.33156464479052"; "M"="18.19273374030685"; "N"="21.11271908439058" }
    4 = @{ "B"="21.90513216058519"; "C"="8.755487409366026"; "D"="7.864124422259168"; "F"="40.63551613417289"; "G"="3.713459174138749"; "L"="10.34279304929871"; "M"="18.15588074907541"; "N"="21.14648747776306" }
    5 = @{ "B"="21.80015319799513"; "C"="8.587323182169532"; "D"="7.870353893361699"; "F"="40.55853134252719"; "G"="3.714657899302115"; "L"="10.34771236352925"; "M"="18.14194507421912"; "N"="21.16073876447976" }
    6 = @{ "B"="21.78278384093613"; "C"="8.559058485035955"; "D"="7.871404014013018"; "F"="40.545916201538"; "G"="3.714859069084084"; "L"="10.34854996725675"; "M"="18.13969667588292"; "N"="21.16313478719852" }
    7 = @{ "B"="21.90371229891824"; "C"="8.753242364591673"; "D"="7.864207380762059"; "F"="40.63446665712617"; "G"="3.713475198374748"; "L"="10.3428580014642"; "M"="18.15568841446002"; "N"="21.1466776905016" }
    8 = @{ "B"="22.44532283022811"; "C"="9.560496249969757"; "D"="7.834739330879237"; "F"="41.04923012833949"; "G"="3.707675986713997"; "L"="10.32075752954645"; "M"="18.23573274125779"; "N"="21.07836242204757" }
    9 = @{ "B"="23.52588843953402"; "C"="10.96506816656416"; "D"="7.785218261366807"; "F"="41.94436233221283"; "G"="3.697399309191982"; "L"="10.28826119426183"; "M"="18.42608625560391"; "N"="20.95990034227226" }
    10 = @{ "B"="24.32211455477156"; "C"="11.88857130425552"; "D"="7.753910341463641"; "F"="42.64738615030132"; "G"="3.69050775640237"; "L"="10.27099514995182"; "M"="18.58558553511013"; "N"="20.88234200254692" }
    11 = @{ "B"="24.6830402169081"; "C"="12.28528165010505"; "D"="7.740777791035937"; "F"="42.97602127330254"; "G"="3.687513683816277"; "L"="10.26457479728427"; "M"="18.66223661839261"; "N"="20.84912213925351" }
    12 = @{ "B"="24.81938664341913"; "C"="12.43215340710257"; "D"="7.735965123093662"; "F"="43.10164391457549"; "G"="3.686400020210465"; "L"="10.26234967191525"; "M"="18.69183456542808"; "N"="20.83683975616675" }
    13 = @{ "B"="24.79003881863762"; "C"="12.40067092779559"; "D"="7.73699447312048"; "F"="43.07453805566423"; "G"="3.686638974659761"; "L"="10.26281972692743"; "M"="18.68543496422908"; "N"="20.83947175373881" }
    14 = @{ "B"="24.69426489841621"; "C"="12.2974320289496"; "D"="7.740378631829378"; "F"="42.98633330524325"; "G"="3.687421659466245"; "L"="10.26438760451779"; "M"="18.66466029887992"; "N"="20.84810569815182" }
    15 = @{ "B"="24.6355537011601"; "C"="12.23375878821958"; "D"="7.742472432010686"; "F"="42.93245563993664"; "G"="3.687903693486038"; "L"="10.26537481479127"; "M"="18.65200916972582"; "N"="20.85343297376583" }
    16 = @{ "B"="24.29848895620395"; "C"="11.86217542556909"; "D"="7.754790978927232"; "F"="42.62607873309545"; "G"="3.690706250077235"; "L"="10.27144358478924"; "M"="18.58065717006498"; "N"="20.88455455277079" }
    17 = @{ "B"="24.09127775906004"; "C"="11.6282357851426"; "D"="7.76263275186186"; "F"="42.44032120733729"; "G"="3.692461524333403"; "L"="10.27553381913148"; "M"="18.537921934259"; "N"="20.90417524998783" }
    18 = @{ "B"="23.97198292434051"; "C"="11.49147702972715"; "D"="7.767247522682351"; "F"="42.33431491094552"; "G"="3.693484384879515"; "L"="10.27802140161491"; "M"="18.51372818662896"; "N"="20.91565457527368" }
    19 = @{ "B"="23.931576933688"; "C"="11.44479422750821"; "D"="7.76882791027989"; "F"="42.29856948433207"; "G"="3.693832991803944"; "L"="10.27888684089517"; "M"="18.5056034665188"; "N"="20.91957457763304" }
    20 = @{ "B"="24.113348446697"; "C"="11.65336697588852"; "D"="7.761787172100671"; "F"="42.4600095067474"; "G"="3.692273299671039"; "L"="10.27508443654619"; "M"="18.54243129810773"; "N"="20.90206650757067" }
    21 = @{ "B"="24.72240604239531"; "C"="12.3278467279709"; "D"="7.739380264307585"; "F"="43.01221002783038"; "G"="3.687191220837391"; "L"="10.26392148764772"; "M"="18.67074694138241"; "N"="20.84556162553007" }
    22 = @{ "B"="25.11848521554799"; "C"="12.74911528276942"; "D"="7.725671042827692"; "F"="43.37991281824877"; "G"="3.683987038138609"; "L"="10.25782723115166"; "M"="18.75793327276437"; "N"="20.810365453521" }
    23 = @{ "B"="24.90731674442152"; "C"="12.52606069794442"; "D"="7.732902094857823"; "F"="43.18307143002205"; "G"="3.685686488181684"; "L"="10.26096996350785"; "M"="18.71110190996131"; "N"="20.82899147458822" }
    24 = @{ "B"="24.10337079597934"; "C"="11.64201221637317"; "D"="7.762169127261211"; "F"="42.45110596451474"; "G"="3.69235835326416"; "L"="10.27528717863393"; "M"="18.54039144444146"; "N"="20.90301924984543" }
    25 = @{ "B"="23.23252584020252"; "C"="10.60417344271383"; "D"="7.797726458841512"; "F"="41.69389572279281"; "G"="3.700063083697186"; "L"="10.29589152043482"; "M"="18.37108442177565"; "N"="20.99028531402024" }
}

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = [double]$rowData[$col]
    }
}
